# Tidied the dataset for 1962
#
# - Insert a new "DEPARTAMENTO" column (B), shifting the rest of the table
#   one column to the right.
# - Re-order the "partidos" / "coalicion" columns (partidos now comes right
#   after "vuelta", coalicion moves next to "Partido Lider").
# - Drop the now-redundant "partidos" values that used to read "total" /
#   "Vincho" / "PLD" for the aggregate rows, and encode the winning party
#   ("PLD" / "FP") directly under "Partido Lider" instead.
# - Add a documentation note about the atomic/yellow columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("B1").Value = "DEPARTAMENTO"
$ws.Range("C1").Value = "Provincia"
$ws.Range("D1").Value = "Municipio"
$ws.Range("E1").Value = "Nivel"
$ws.Range("F1").Value = "vuelta"
$ws.Range("G1").Value = "partidos"
$ws.Range("H1").Value = "coalicion"
$ws.Range("I1").Value = "Partido Lider"
$ws.Range("J1").Value = "siglas"
$ws.Range("K1").Value = "# de votos"
$ws.Range("L1").Value = "arrastre"
$ws.Range("M1").Value = "# de colegios"
$ws.Range("N1").Value = "# inscritos"
$ws.Range("O1").Value = "votos validos"
$ws.Range("P1").Value = "votos nulos"

# --- Row 2: 2020, vuelta 1, PLD y aliados ---------------------------------
$ws.Range("C2").Value = "Santiago"
$ws.Range("D2").Value = "san jose de las matas"
$ws.Range("E2").Value = "generales"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "Partido de la Liberacion Dominicana"
$ws.Range("H2").Value = "PLD y aliados"
$ws.Range("I2").Value = "PLD"
$ws.Range("J2").Value = "PLD"
$ws.Range("K2").Value = "NA"
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = "Valido solo desde 1998"

# --- Row 3: 2020, vuelta 2, PLD y aliados ---------------------------------
$ws.Range("C3").Value = "Santiago"
$ws.Range("D3").Value = "san jose de las matas"
$ws.Range("E3").Value = "generales"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = "Vincho"
$ws.Range("H3").Value = "PLD y aliados"
$ws.Range("I3").Value = "PLD"
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = "NA"
$ws.Range("L3").Value = 2

# --- Row 4: 2020, vuelta 1, FP y aliados ----------------------------------
$ws.Range("C4").Value = "Santiago"
$ws.Range("D4").Value = "san jose de las matas"
$ws.Range("E4").Value = "generales"
$ws.Range("F4").Value = 1
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = "FP y aliados"
$ws.Range("I4").Value = "FP"
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 5

# --- Row 5: 2020, vuelta 2, FP y aliados ----------------------------------
$ws.Range("C5").Value = "Santiago"
$ws.Range("D5").Value = "san jose de las matas"
$ws.Range("E5").Value = "generales"
$ws.Range("F5").Value = 2
$ws.Range("G5").ClearContents()
$ws.Range("H5").Value = "FP y aliados"
$ws.Range("I5").Value = "FP"
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 700

# --- Row 6: 2024, vuelta 1, PLD y aliados ---------------------------------
$ws.Range("C6").Value = "Santiago"
$ws.Range("D6").Value = "san jose de las matas"
$ws.Range("E6").Value = "generales"
$ws.Range("F6").Value = 1
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = "PLD y aliados"
$ws.Range("I6").Value = "PLD"
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = 300

# --- Row 7: 2024, vuelta 2, PLD y aliados ---------------------------------
$ws.Range("C7").Value = "Santiago"
$ws.Range("D7").Value = "san jose de las matas"
$ws.Range("E7").Value = "generales"
$ws.Range("F7").Value = 2
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = "PLD y aliados"
$ws.Range("I7").Value = "PLD"
$ws.Range("J7").ClearContents()
$ws.Range("K7").Value = 1000

# Columns B (Provincia) through B (old) held "Santiago" for every row; that
# now lives in column C, so clear out the stale leftovers in B.
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()

# --- Documentation note ----------------------------------------------------
$ws.Range("A11").Value = "Observacion: aquellas columnas en amarillo definen una observacion atomica"

# --- Selection / view state, to mirror the saved workbook -----------------
$ws.Range("H1").Select()
